$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.763.57'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '1.884.22'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '237.46'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '0.4781'
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('D8').Value = '0.2842'
$ws.Range('E8').Value = '  +4.49%  '
$ws.Range('D9').Value = '0.06488'
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('D10').Value = '18.87'
$ws.Range('E10').Value = '  +17.31%  '
$ws.Range('D11').Value = '1.891.90'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').Value = '0.07555'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '94.10'
$ws.Range('E13').Value = '  +12.42%  '
$ws.Range('D14').Value = '5.098'
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('D15').Value = '0.6501'
$ws.Range('E15').Value = '  +4.90%  '
$ws.Range('D16').Value = '295.87'
$ws.Range('E16').Value = '  +30.94%  '
$ws.Range('D17').Value = '30.758.78'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').Value = '13.14'
$ws.Range('E18').Value = '  +6.63%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '0.000007475'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Value = '2.144.62'
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '5.136'
$ws.Range('E23').Value = '  +5.09%  '
$ws.Range('D24').Value = '6.130'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '167.98'
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.281'
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('D27').Value = '19.52'
$ws.Range('E27').Value = '  +10.06%  '
$ws.Range('D28').Value = '1.947'
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').Value = '0.1059'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('D30').Value = '1.356'
$ws.Range('E30').Value = '  -1.34%  '
$ws.Range('D31').Value = '4.178'
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').Value = '3.947'
$ws.Range('E32').Value = '  +3.85%  '
$ws.Range('D33').Value = '0.05019'
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('D34').Value = '1.167'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('D35').Value = '0.7200'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('D36').Value = '2.716'
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('D37').Value = '0.01943'
$ws.Range('E37').Value = '  +4.39%  '
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('D39').Value = '2.055'
$ws.Range('E39').Value = '  +7.49%  '
$ws.Range('D40').Value = '0.8955'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').Value = '106.96'
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '0.4189'
$ws.Range('E43').Value = '  +4.48%  '
$ws.Range('D44').Value = '5.571'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').Value = '65.94'
$ws.Range('E45').Value = '  +10.43%  '
$ws.Range('D46').Value = '7.321'
$ws.Range('E46').Value = '  +4.13%  '
$ws.Range('D47').Value = '0.1221'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '8.890'
$ws.Range('E48').Value = '  +3.53%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '34.56'
$ws.Range('E49').Value = '  +5.11%  '
$ws.Range('D50').Value = '0.05633'
$ws.Range('E50').Value = '  +2.24%  '
$ws.Range('D51').Value = '1.382'
$ws.Range('E51').Value = '  +1.97%  '
